# Update column C values on Sheet1 (RandomForest imputed values)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value = -10.83959999999999
$ws.Range("C10").Value = -12.8093
$ws.Range("C12").Value = -14.35450000000001
$ws.Range("C18").Value = -14.22130000000001
